# Notifications.xlsx: flip Runmode column from "N" to "Y" for the
# Notifications test cases (rows 4-26 on the "Test Cases" sheet) and
# update the sheet's active selection to D3:D26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

for ($r = 4; $r -le 26; $r++) {
    $ws.Cells.Item($r, 4).Value = "Y"
}

$ws.Activate() | Out-Null
$ws.Range("D3:D26").Select() | Out-Null
